$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.217416695739445
$ws.Cells.Item(2, 3).Value = 0.03324176402566081
$ws.Cells.Item(2, 4).Value = 0.4669838078361579
$ws.Cells.Item(2, 5).Value = 0.1325495583738903
$ws.Cells.Item(2, 7).Value = 2.558819655718906
$ws.Cells.Item(2, 8).Value = 2.076381815815182
$ws.Cells.Item(2, 9).Value = 2.071961212152104
$ws.Cells.Item(2, 10).Value = 0.04081147924233086
$ws.Cells.Item(2, 11).Value = 0.9240908947411128
$ws.Cells.Item(2, 12).Value = 0.4890084694475476
$ws.Cells.Item(2, 13).Value = 0.3818061174113829
$ws.Cells.Item(2, 14).Value = 3.561060411187135
$ws.Cells.Item(3, 2).Value = 1.187814966634278
$ws.Cells.Item(3, 3).Value = 0.02898488964095236
$ws.Cells.Item(3, 4).Value = 0.4657646514010594
$ws.Cells.Item(3, 5).Value = 0.1329194048748672
$ws.Cells.Item(3, 7).Value = 2.561295778591813
$ws.Cells.Item(3, 8).Value = 2.082958715696492
$ws.Cells.Item(3, 9).Value = 2.078489777806951
$ws.Cells.Item(3, 10).Value = 0.04050603993829505
$ws.Cells.Item(3, 11).Value = 0.8921205157163286
$ws.Cells.Item(3, 12).Value = 0.4876249809881088
$ws.Cells.Item(3, 13).Value = 0.376596550606525
$ws.Cells.Item(3, 14).Value = 3.583554135965265
$ws.Cells.Item(4, 2).Value = 1.170239634185577
$ws.Cells.Item(4, 3).Value = 0.02637132571686607
$ws.Cells.Item(4, 4).Value = 0.4651943940453265
$ws.Cells.Item(4, 5).Value = 0.1331792236262235
$ws.Cells.Item(4, 7).Value = 2.563847844849022
$ws.Cells.Item(4, 8).Value = 2.087672364156973
$ws.Cells.Item(4, 9).Value = 2.083210593413831
$ws.Cells.Item(4, 10).Value = 0.04031352725831816
$ws.Cells.Item(4, 11).Value = 0.8729720700697214
$ws.Cells.Item(4, 12).Value = 0.4869615784810932
$ws.Cells.Item(4, 13).Value = 0.3735723077943192
$ws.Cells.Item(4, 14).Value = 3.598206726410702
$ws.Cells.Item(5, 2).Value = 1.163228973486213
$ws.Cells.Item(5, 3).Value = 0.02530629586736666
$ws.Cells.Item(5, 4).Value = 0.4650069426550516
$ws.Cells.Item(5, 5).Value = 0.1332933528090088
$ws.Cells.Item(5, 7).Value = 2.565147399567834
$ws.Cells.Item(5, 8).Value = 2.089763244058545
$ws.Cells.Item(5, 9).Value = 2.085313618574268
$ws.Cells.Item(5, 10).Value = 0.04023382162943712
$ws.Cells.Item(5, 11).Value = 0.8652903515660739
$ws.Cells.Item(5, 12).Value = 0.486738123354499
$ws.Cells.Item(5, 13).Value = 0.3723838908492851
$ws.Cells.Item(5, 14).Value = 3.604389424445976
$ws.Cells.Item(6, 2).Value = 1.162074019167989
$ws.Cells.Item(6, 3).Value = 0.02512944887114088
$ws.Cells.Item(6, 4).Value = 0.4649785331687468
$ws.Cells.Item(6, 5).Value = 0.1333128028076835
$ws.Cells.Item(6, 7).Value = 2.565378870838458
$ws.Cells.Item(6, 8).Value = 2.090120708189687
$ws.Cells.Item(6, 9).Value = 2.085673655377597
$ws.Cells.Item(6, 10).Value = 0.04022051061696885
$ws.Cells.Item(6, 11).Value = 0.8640221524822493
$ws.Cells.Item(6, 12).Value = 0.4867038535764721
$ws.Cells.Item(6, 13).Value = 0.3721892153189756
$ws.Cells.Item(6, 14).Value = 3.605428842297329
$ws.Cells.Item(7, 2).Value = 1.170144472176105
$ws.Cells.Item(7, 3).Value = 0.02635696230832707
$ws.Cells.Item(7, 4).Value = 0.4651916839499251
$ws.Cells.Item(7, 5).Value = 0.1331807293774201
$ws.Cells.Item(7, 7).Value = 2.563864319948351
$ws.Cells.Item(7, 8).Value = 2.087699873797746
$ws.Cells.Item(7, 9).Value = 2.083238229552293
$ws.Cells.Item(7, 10).Value = 0.04031245740873324
$ws.Cells.Item(7, 11).Value = 0.8728679795994196
$ws.Cells.Item(7, 12).Value = 0.4869583749120139
$ws.Cells.Item(7, 13).Value = 0.3735561021296014
$ws.Cells.Item(7, 14).Value = 3.59828925137888
$ws.Cells.Item(8, 2).Value = 1.207085720597519
$ws.Cells.Item(8, 3).Value = 0.03177394123063948
$ws.Cells.Item(8, 4).Value = 0.4665264905106028
$ws.Cells.Item(8, 5).Value = 0.1326703001261382
$ws.Cells.Item(8, 7).Value = 2.559459364686603
$ws.Cells.Item(8, 8).Value = 2.078509464551274
$ws.Cells.Item(8, 9).Value = 2.074064523829321
$ws.Cells.Item(8, 10).Value = 0.04070719260649369
$ws.Cells.Item(8, 11).Value = 0.9129677722209237
$ws.Cells.Item(8, 12).Value = 0.4884928827995765
$ws.Cells.Item(8, 13).Value = 0.3799737200066069
$ws.Cells.Item(8, 14).Value = 3.568641676159189
$ws.Cells.Item(9, 2).Value = 1.284273594322713
$ws.Cells.Item(9, 3).Value = 0.04239964891763748
$ws.Cells.Item(9, 4).Value = 0.4705552283794106
$ws.Cells.Item(9, 5).Value = 0.1319281999602762
$ws.Cells.Item(9, 7).Value = 2.559003942157162
$ws.Cells.Item(9, 8).Value = 2.065838518598511
$ws.Cells.Item(9, 9).Value = 2.061720553001365
$ws.Cells.Item(9, 10).Value = 0.04144206394597028
$ws.Cells.Item(9, 11).Value = 0.9954134732071793
$ws.Cells.Item(9, 12).Value = 0.4929746593599162
$ws.Cells.Item(9, 13).Value = 0.3939386880764957
$ws.Cells.Item(9, 14).Value = 3.517174392644428
$ws.Cells.Item(10, 2).Value = 1.343862327045315
$ws.Cells.Item(10, 3).Value = 0.05021122116740173
$ws.Cells.Item(10, 4).Value = 0.4743714356435902
$ws.Cells.Item(10, 5).Value = 0.1315396742185868
$ws.Cells.Item(10, 7).Value = 2.563656445908606
$ws.Cells.Item(10, 8).Value = 2.059782621329134
$ws.Cells.Item(10, 9).Value = 2.056086925932227
$ws.Cells.Item(10, 10).Value = 0.04195844385388803
$ws.Cells.Item(10, 11).Value = 1.058302828189795
$ws.Cells.Item(10, 12).Value = 0.4971610773646944
$ws.Cells.Item(10, 13).Value = 0.4050360488724181
$ws.Cells.Item(10, 14).Value = 3.483422001186106
$ws.Cells.Item(11, 2).Value = 1.371593297010435
$ws.Cells.Item(11, 3).Value = 0.0537666555054841
$ws.Cells.Item(11, 4).Value = 0.4762926777741114
$ws.Cells.Item(11, 5).Value = 0.1313967192413674
$ws.Cells.Item(11, 7).Value = 2.566855913310718
$ws.Cells.Item(11, 8).Value = 2.057732334932325
$ws.Cells.Item(11, 9).Value = 2.054268848389022
$ws.Cells.Item(11, 10).Value = 0.04218832961864827
$ws.Cells.Item(11, 11).Value = 1.087415004264528
$ws.Cells.Item(11, 12).Value = 0.4992588394517981
$ws.Cells.Item(11, 13).Value = 0.4102655955871981
$ws.Cells.Item(11, 14).Value = 3.468947549965165
$ws.Cells.Item(12, 2).Value = 1.382183610912449
$ws.Cells.Item(12, 3).Value = 0.05511331948071074
$ws.Cells.Item(12, 4).Value = 0.477046756744258
$ws.Cells.Item(12, 5).Value = 0.1313474251230495
$ws.Cells.Item(12, 7).Value = 2.568223151999092
$ws.Cells.Item(12, 8).Value = 2.057057102411079
$ws.Cells.Item(12, 9).Value = 2.053687360767221
$ws.Cells.Item(12, 10).Value = 0.04227466541925651
$ws.Cells.Item(12, 11).Value = 1.098511193924651
$ws.Cells.Item(12, 12).Value = 0.5000809260103551
$ws.Cells.Item(12, 13).Value = 0.4122718653525013
$ws.Cells.Item(12, 14).Value = 3.463592884005273
$ws.Cells.Item(13, 2).Value = 1.379898838160329
$ws.Cells.Item(13, 3).Value = 0.05482327766654294
$ws.Cells.Item(13, 4).Value = 0.4768831731159224
$ws.Cells.Item(13, 5).Value = 0.1313578265381103
$ws.Cells.Item(13, 7).Value = 2.567921770269265
$ws.Cells.Item(13, 8).Value = 2.057198028777918
$ws.Cells.Item(13, 9).Value = 2.053807838286687
$ws.Cells.Item(13, 10).Value = 0.04225610327285345
$ws.Cells.Item(13, 11).Value = 1.096118232668459
$ws.Cells.Item(13, 12).Value = 0.4999026440361405
$ws.Cells.Item(13, 13).Value = 0.4118386268125178
$ws.Cells.Item(13, 14).Value = 3.464740481241947
$ws.Cells.Item(14, 2).Value = 1.372462783705657
$ws.Cells.Item(14, 3).Value = 0.05387744025146901
$ws.Cells.Item(14, 4).Value = 0.4763541847904094
$ws.Cells.Item(14, 5).Value = 0.1313925668942151
$ws.Cells.Item(14, 7).Value = 2.566965277381485
$ws.Cells.Item(14, 8).Value = 2.057674756472437
$ws.Cells.Item(14, 9).Value = 2.054218865638362
$ws.Cells.Item(14, 10).Value = 0.04219544685933485
$ws.Cells.Item(14, 11).Value = 1.088326452733924
$ws.Cells.Item(14, 12).Value = 0.4993259181938612
$ws.Cells.Item(14, 13).Value = 0.4104301331216504
$ws.Cells.Item(14, 14).Value = 3.468504483234845
$ws.Cells.Item(15, 2).Value = 1.367919587027274
$ws.Cells.Item(15, 3).Value = 0.05329812704638925
$ws.Cells.Item(15, 4).Value = 0.4760336186212299
$ws.Cells.Item(15, 5).Value = 0.1314144760994651
$ws.Cells.Item(15, 7).Value = 2.566399668723875
$ws.Cells.Item(15, 8).Value = 2.057979936388449
$ws.Cells.Item(15, 9).Value = 2.054484560446447
$ws.Cells.Item(15, 10).Value = 0.04215819982404589
$ws.Cells.Item(15, 11).Value = 1.083563133625603
$ws.Cells.Item(15, 12).Value = 0.4989762628355123
$ws.Cells.Item(15, 13).Value = 0.4095707665681729
$ws.Cells.Item(15, 14).Value = 3.470826516459951
$ws.Cells.Item(16, 2).Value = 1.342062551905968
$ws.Cells.Item(16, 3).Value = 0.04997890562300711
$ws.Cells.Item(16, 4).Value = 0.4742495975883543
$ws.Cells.Item(16, 5).Value = 0.1315496949447805
$ws.Cells.Item(16, 7).Value = 2.563469140526053
$ws.Cells.Item(16, 8).Value = 2.059930777269614
$ws.Cells.Item(16, 9).Value = 2.056220717463709
$ws.Cells.Item(16, 10).Value = 0.04194331969958576
$ws.Cells.Item(16, 11).Value = 1.056410379673622
$ws.Cells.Item(16, 12).Value = 0.4970278662534042
$ws.Cells.Item(16, 13).Value = 0.4046979239326234
$ws.Cells.Item(16, 14).Value = 3.48438564174247
$ws.Cells.Item(17, 2).Value = 1.326359502488941
$ws.Cells.Item(17, 3).Value = 0.04794317685986016
$ws.Cells.Item(17, 4).Value = 0.4732025364909447
$ws.Cells.Item(17, 5).Value = 0.1316412876967057
$ws.Cells.Item(17, 7).Value = 2.561948695279625
$ws.Cells.Item(17, 8).Value = 2.061307905198248
$ws.Cells.Item(17, 9).Value = 2.057476449692118
$ws.Cells.Item(17, 10).Value = 0.04181021576904698
$ws.Cells.Item(17, 11).Value = 1.039881743803022
$ws.Cells.Item(17, 12).Value = 0.4958820410335676
$ws.Cells.Item(17, 13).Value = 0.4017549450763838
$ws.Cells.Item(17, 14).Value = 3.492929030326202
$ws.Cells.Item(18, 2).Value = 1.317386264374591
$ws.Cells.Item(18, 3).Value = 0.04677245799108221
$ws.Cells.Item(18, 4).Value = 0.4726177347283169
$ws.Cells.Item(18, 5).Value = 0.1316971507626779
$ws.Cells.Item(18, 7).Value = 2.561176118207783
$ws.Cells.Item(18, 8).Value = 2.062166319197075
$ws.Cells.Item(18, 9).Value = 2.058268813398371
$ws.Cells.Item(18, 10).Value = 0.04173318531100101
$ws.Cells.Item(18, 11).Value = 1.030422322287052
$ws.Cells.Item(18, 12).Value = 0.4952411961974263
$ws.Cells.Item(18, 13).Value = 0.4000792934096893
$ws.Cells.Item(18, 14).Value = 3.497925761074264
$ws.Cells.Item(19, 2).Value = 1.314358184086274
$ws.Cells.Item(19, 3).Value = 0.04637610263843328
$ws.Cells.Item(19, 4).Value = 0.4724227289890877
$ws.Cells.Item(19, 5).Value = 0.1317166120076987
$ws.Cells.Item(19, 7).Value = 2.560932047757404
$ws.Cells.Item(19, 8).Value = 2.062468359672437
$ws.Cells.Item(19, 9).Value = 2.058549137110163
$ws.Cells.Item(19, 10).Value = 0.04170702281037819
$ws.Cells.Item(19, 11).Value = 1.027227678898413
$ws.Cells.Item(19, 12).Value = 0.4950273464440329
$ws.Cells.Item(19, 13).Value = 0.3995148822422365
$ws.Cells.Item(19, 14).Value = 3.499631789312055
$ws.Cells.Item(20, 2).Value = 1.328025043089923
$ws.Cells.Item(20, 3).Value = 0.04815986502147496
$ws.Cells.Item(20, 4).Value = 0.473312193599611
$ws.Cells.Item(20, 5).Value = 0.1316312083716635
$ws.Cells.Item(20, 7).Value = 2.562099999061161
$ws.Cells.Item(20, 8).Value = 2.06115444413021
$ws.Cells.Item(20, 9).Value = 2.057335520563747
$ws.Cells.Item(20, 10).Value = 0.04182443378830669
$ws.Cells.Item(20, 11).Value = 1.041636340895906
$ws.Cells.Item(20, 12).Value = 0.4960021327594291
$ws.Cells.Item(20, 13).Value = 0.402066464168108
$ws.Cells.Item(20, 14).Value = 3.492011002599483
$ws.Cells.Item(21, 2).Value = 1.374644516319819
$ws.Cells.Item(21, 3).Value = 0.05415524734932831
$ws.Cells.Item(21, 4).Value = 0.4765088417385073
$ws.Cells.Item(21, 5).Value = 0.1313822316067892
$ws.Cells.Item(21, 7).Value = 2.567241998216957
$ws.Cells.Item(21, 8).Value = 2.057531985549957
$ws.Cells.Item(21, 9).Value = 2.054095234402652
$ws.Cells.Item(21, 10).Value = 0.04221328253903422
$ws.Cells.Item(21, 11).Value = 1.090613135094173
$ws.Cells.Item(21, 12).Value = 0.499494565257379
$ws.Cells.Item(21, 13).Value = 0.4108431384093976
$ws.Cells.Item(21, 14).Value = 3.467395471581789
$ws.Cells.Item(22, 2).Value = 1.405632717996411
$ws.Cells.Item(22, 3).Value = 0.05807532498714352
$ws.Cells.Item(22, 4).Value = 0.4787527267832132
$ws.Cells.Item(22, 5).Value = 0.1312477125061129
$ws.Cells.Item(22, 7).Value = 2.571509925351847
$ws.Cells.Item(22, 8).Value = 2.055754106749504
$ws.Cells.Item(22, 9).Value = 2.052601026972113
$ws.Cells.Item(22, 10).Value = 0.04246324088550857
$ws.Cells.Item(22, 11).Value = 1.123042041133317
$ws.Cells.Item(22, 12).Value = 0.5019385419448383
$ws.Cells.Item(22, 13).Value = 0.4167304277383579
$ws.Cells.Item(22, 14).Value = 3.452045107954916
$ws.Cells.Item(23, 2).Value = 1.389046344819121
$ws.Cells.Item(23, 3).Value = 0.05598293805520882
$ws.Cells.Item(23, 4).Value = 0.4775409983963357
$ws.Cells.Item(23, 5).Value = 0.1313169335975051
$ws.Cells.Item(23, 7).Value = 2.569149052803908
$ws.Cells.Item(23, 8).Value = 2.056649094073265
$ws.Cells.Item(23, 9).Value = 2.053341497294461
$ws.Cells.Item(23, 10).Value = 0.04233021406101756
$ws.Cells.Item(23, 11).Value = 1.105695836506072
$ws.Cells.Item(23, 12).Value = 0.5006194014527523
$ws.Cells.Item(23, 13).Value = 0.4135744720357337
$ws.Cells.Item(23, 14).Value = 3.460170421916615
$ws.Cells.Item(24, 2).Value = 1.327271881903329
$ws.Cells.Item(24, 3).Value = 0.04806190140094202
$ws.Cells.Item(24, 4).Value = 0.4732625641398869
$ws.Cells.Item(24, 5).Value = 0.1316357552494605
$ws.Cells.Item(24, 7).Value = 2.562031278295592
$ws.Cells.Item(24, 8).Value = 2.061223616139046
$ws.Cells.Item(24, 9).Value = 2.057399015232761
$ws.Cells.Item(24, 10).Value = 0.04181800740221675
$ws.Cells.Item(24, 11).Value = 1.040842953235028
$ws.Cells.Item(24, 12).Value = 0.4959477835100046
$ws.Cells.Item(24, 13).Value = 0.4019255755860556
$ws.Cells.Item(24, 14).Value = 3.492425778076814
$ws.Cells.Item(25, 2).Value = 1.262885633515481
$ws.Cells.Item(25, 3).Value = 0.03952447715057872
$ws.Cells.Item(25, 4).Value = 0.4693146863808266
$ws.Cells.Item(25, 5).Value = 0.1321013608738664
$ws.Cells.Item(25, 7).Value = 2.558251275891607
$ws.Cells.Item(25, 8).Value = 2.068694369543252
$ws.Cells.Item(25, 9).Value = 2.064456173008146
$ws.Cells.Item(25, 10).Value = 0.04124742050221641
$ws.Cells.Item(25, 11).Value = 0.9727024783106799
$ws.Cells.Item(25, 12).Value = 0.4916049914320126
$ws.Cells.Item(25, 13).Value = 0.3900134680676786
$ws.Cells.Item(25, 14).Value = 3.530384113102848
